$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 - sample 29
$ws.Range("C29").Copy($ws.Range("C30"))
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "social"
$ws.Range("C30").Value = 44380
$ws.Range("D30").Value = "MCAST"
$ws.Range("E30").Value = "lookalike"
$ws.Range("F30").Value = "opportunity"
$ws.Range("G30").Value = "mt"
$ws.Range("H30").Value = "no"
$ws.Range("I30").Value = "click to see mature content"

# Row 31 - sample 30 (delivery message)
$ws.Range("C29").Copy($ws.Range("C31"))
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "msg"
$ws.Range("C31").Value = 44380
$ws.Range("D31").Value = "MCAST"
$ws.Range("E31").Value = "lookalike"
$ws.Range("F31").Value = "delivery"
$ws.Range("G31").Value = "mt"
$ws.Range("H31").Value = "no"
$ws.Range("I31").Value = "click to confirm delivery"

# Row 32 - sample 30 (second entry, account blocked)
$ws.Range("C29").Copy($ws.Range("C32"))
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "msg"
$ws.Range("C32").Value = 44380
$ws.Range("D32").Value = "MCAST"
$ws.Range("E32").Value = "shortened"
$ws.Range("F32").Value = "lockout"
$ws.Range("G32").Value = "mt"
$ws.Range("H32").Value = "no"
$ws.Range("I32").Value = "account blocked"

[void]$ws.Range("I32").Select()
